# This workbook is the "phieu nhap" (goods-receipt) import template.
# The sheet was previously mislabeled as the "export production" template;
# rename it to reflect that it is actually the warehouse import ("nhap kho")
# template, per the commit "Update import templates ...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Template nhập kho"
